# BDD.xlsx edit: add "Type / Marque / Modele / Annee" columns (C:F) in
# front of the existing Nom/Voyou/Descriptif block, fill them in per row,
# and tidy up the selection / formatting of the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 4 new columns at C (old C:E -> G:I) shifting everything right
# ---------------------------------------------------------------------
$ws.Range("C1:F1").EntireColumn.Insert(-4161)

# Give the 4 new columns the same width as the rest of the table
$ws.Range("C1:F1").EntireColumn.ColumnWidth = 19.33

# ---------------------------------------------------------------------
# 2. Header row (row 2): Type / Marque / Modele / Annee
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "Type"
$ws.Range("D2").Value = "Marque"
$ws.Range("E2").Value = "Modele"
$ws.Range("F2").Value = "Annee"

# Match the header look of the rest of row 2 (grey fill, centered, medium
# top border, thin separators)
$hdr = $ws.Range("C2:F2")
$hdr.Interior.Pattern = 1
$hdr.Interior.PatternColorIndex = -4105
$hdr.Interior.ThemeColor = 1
$hdr.Interior.TintAndShade = -0.249977111117893
$hdr.HorizontalAlignment = -4108
$hdr.Borders.Item(7).LineStyle = 1
$hdr.Borders.Item(7).Weight = 2
$hdr.Borders.Item(7).ColorIndex = 1
$hdr.Borders.Item(10).LineStyle = 1
$hdr.Borders.Item(10).Weight = -4138
$hdr.Borders.Item(10).ColorIndex = 1

# ---------------------------------------------------------------------
# 3. Body rows 3-27: Type is always "Voiture", Modele is always "Choix";
#    Marque and Annee vary per row (per the BDD data)
# ---------------------------------------------------------------------
$ws.Range("C3:C27").Value = "Voiture"
$ws.Range("E3:E27").Value = "Choix"

$cars = @(
    @{ Row = 3;  Marque = "BMW";      Annee = 2015 },
    @{ Row = 4;  Marque = "Audi";     Annee = 2020 },
    @{ Row = 5;  Marque = "Mercedes"; Annee = 2015 },
    @{ Row = 6;  Marque = "Audi";     Annee = 2013 },
    @{ Row = 7;  Marque = "Mercedes"; Annee = 2019 },
    @{ Row = 8;  Marque = "Renault";  Annee = 2018 },
    @{ Row = 9;  Marque = "BMW";      Annee = 2014 },
    @{ Row = 10; Marque = "Ford";     Annee = 2016 },
    @{ Row = 11; Marque = "Peugeot";  Annee = 2022 },
    @{ Row = 12; Marque = "Audi";     Annee = 2020 },
    @{ Row = 13; Marque = "Mercedes"; Annee = 2019 },
    @{ Row = 14; Marque = "Ford";     Annee = 2023 },
    @{ Row = 15; Marque = "Citroen";  Annee = 2022 },
    @{ Row = 16; Marque = "Ford";     Annee = 2013 },
    @{ Row = 17; Marque = "Citroen";  Annee = 2012 },
    @{ Row = 18; Marque = "Citroen";  Annee = 2015 },
    @{ Row = 19; Marque = "Mercedes"; Annee = 2022 },
    @{ Row = 20; Marque = "BMW";      Annee = 2018 },
    @{ Row = 21; Marque = "Citroen";  Annee = 2010 },
    @{ Row = 22; Marque = "Mercedes"; Annee = 2013 },
    @{ Row = 23; Marque = "Audi";     Annee = 2022 },
    @{ Row = 24; Marque = "BMW";      Annee = 2019 },
    @{ Row = 25; Marque = "Mercedes"; Annee = 2022 },
    @{ Row = 26; Marque = "Audi";     Annee = 2017 },
    @{ Row = 27; Marque = "Citroen";  Annee = 2015 }
)

foreach ($car in $cars) {
    $ws.Range("D" + $car.Row).Value = $car.Marque
    $ws.Range("F" + $car.Row).Value = $car.Annee
}

# ---------------------------------------------------------------------
# 4. Borders / alignment for the new body block (C3:F27)
# ---------------------------------------------------------------------
$body = $ws.Range("C3:F27")
$body.HorizontalAlignment = -4108
$body.Borders.Item(7).LineStyle = 1
$body.Borders.Item(7).Weight = 2
$body.Borders.Item(7).ColorIndex = 1
$body.Borders.Item(9).LineStyle = 1
$body.Borders.Item(9).Weight = 2
$body.Borders.Item(9).ColorIndex = 1
$body.Borders.Item(10).LineStyle = 1
$body.Borders.Item(10).Weight = 2
$body.Borders.Item(10).ColorIndex = 1
$body.Borders.Item(11).LineStyle = 1
$body.Borders.Item(11).Weight = 2
$body.Borders.Item(11).ColorIndex = 1
$body.Borders.Item(12).LineStyle = 1
$body.Borders.Item(12).Weight = 2
$body.Borders.Item(12).ColorIndex = 1

$ws.Range("C3:F26").Borders.Item(8).LineStyle = 1
$ws.Range("C3:F26").Borders.Item(8).Weight = 2
$ws.Range("C3:F26").Borders.Item(8).ColorIndex = 1

$left = $ws.Range("C3:C27")
$left.Borders.Item(7).LineStyle = 1
$left.Borders.Item(7).Weight = -4138
$left.Borders.Item(7).ColorIndex = 1

$right = $ws.Range("F3:F27")
$right.Borders.Item(10).LineStyle = 1
$right.Borders.Item(10).Weight = 2
$right.Borders.Item(10).ColorIndex = 1

$ws.Range("C27:F27").Borders.Item(9).LineStyle = 1
$ws.Range("C27:F27").Borders.Item(9).Weight = 2
$ws.Range("C27:F27").Borders.Item(9).ColorIndex = 1

# ---------------------------------------------------------------------
# 5. Selection, as left by the author after finishing the edit
# ---------------------------------------------------------------------
$ws.Range("L18").Select()
